# Reto 3 - Graficas: "Actualizacion Documento y Excel con el Req 5"
#
# The author filled in the previously-empty "Req 5" column (column G) of
# the results table on sheet "Tiempo de ejecucion " with the measured
# execution times (matching the pattern already present for Req 1..Req 4):
#   G5  = 0        (sample size 0%)
#   G6  = 15.625   (5%)
#   G7  = 15.625   (10%)
#   G8  = 15.625   (20%)
#   G9  = 15.625   (30%)
#   G10 = 15.625   (50%)
#   G11 = 15.625   (80%)
#   G12 = 15.625   (100%)
# and left the selection sitting on G6 afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tiempo de ejecucion ")

$ws.Range("G5").Value  = 0
$ws.Range("G6").Value  = 15.625
$ws.Range("G7").Value  = 15.625
$ws.Range("G8").Value  = 15.625
$ws.Range("G9").Value  = 15.625
$ws.Range("G10").Value = 15.625
$ws.Range("G11").Value = 15.625
$ws.Range("G12").Value = 15.625

# Leave the worksheet selection where the author left it when saving.
$ws.Activate()
$ws.Range("G6").Select()
